$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50 - 四方坪站
$ws.Range("A50").Value = 46016
$ws.Range("B50").Value = "四方坪站"
$ws.Range("C50").Value = 10431.219999999999
$ws.Range("D50").Value = 8802.5300000000007
$ws.Range("E50").Value = 3547.74
$ws.Range("F50").Value = 438

# Row 51 - 高岭站
$ws.Range("A51").Value = 46016
$ws.Range("B51").Value = "高岭站"
$ws.Range("C51").Value = 6540.3
$ws.Range("D51").Value = 5806.07
$ws.Range("E51").Value = 1724.04
$ws.Range("F51").Value = 222

$ws.Range("I51").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
